$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 6-9 are new VIN-upload rows appended below the existing rows 2-5.
# Each new row reuses the same vehicle attributes as its row-2..row-5
# counterpart (same VIN/MAKE/etc.), but with:
#   - column B changed from "SYMBOL_2000" to "SYMBOL_2017"
#   - columns AC..AF changed from the single repeated code to four distinct
#     per-row codes (BI00n / PD00n / UM00n / MP00n)
# The cell-level formatting (alternating Normal/Good/Bad highlighting for
# the ENTRYDATE/VALID columns, left alignment everywhere else) matches the
# corresponding source row exactly, so we copy the formats over afterwards
# instead of re-deriving the style for every single cell.
# ---------------------------------------------------------------------------

$rowData = @(
    @{ Row=6; SrcRow=2; D="TOYOTA"; E="TOYOTA"; F="Gt";          AC="BI001"; AD="PD001"; AE="UM001"; AF="MP001"; AG=20010101; AH="Y" },
    @{ Row=7; SrcRow=3; D="UT_SS";  E="UT_SS";  F="invalidVIN";  AC="BI002"; AD="PD002"; AE="UM002"; AF="MP002"; AG=20000101; AH="N" },
    @{ Row=8; SrcRow=4; D="UT_SS";  E="UT_SS";  F="SecondValid"; AC="BI003"; AD="PD003"; AE="UM003"; AF="MP003"; AG=20150101; AH="Y" },
    @{ Row=9; SrcRow=5; D="UT_SS";  E="UT_SS";  F="ThirdValid";  AC="BI004"; AD="PD004"; AE="UM004"; AF="MP004"; AG=20190101; AH="Y" }
)

foreach ($rd in $rowData) {
    $r = $rd.Row

    $ws.Range("A" + $r).Value  = "XXXKN3DD&E"
    $ws.Range("B" + $r).Value  = "SYMBOL_2017"
    $ws.Range("C" + $r).Value  = 2018
    $ws.Range("D" + $r).Value  = $rd.D
    $ws.Range("E" + $r).Value  = $rd.E
    $ws.Range("F" + $r).Value  = $rd.F
    $ws.Range("G" + $r).Value  = "MDX ADVANCE"
    $ws.Range("H" + $r).Value  = 53080
    $ws.Range("I" + $r).Value  = "WAG"
    $ws.Range("J" + $r).Value  = "UT_SS"
    $ws.Range("K" + $r).Value  = "SUV"
    $ws.Range("L" + $r).Value  = "UT_SS"
    $ws.Range("M" + $r).Value  = "WAG"
    $ws.Range("N" + $r).Value  = "4.5L V10"
    $ws.Range("O" + $r).Value  = 8
    $ws.Range("P" + $r).Value  = "G"
    $ws.Range("Q" + $r).Value  = 214
    $ws.Range("R" + $r).Value  = "2WD"
    $ws.Range("S" + $r).Value  = 2
    $ws.Range("T" + $r).Value  = "000R"
    $ws.Range("U" + $r).Value  = "DUAL AIR BAGS FRONT"
    $ws.Range("V" + $r).Value  = 2
    $ws.Range("W" + $r).Value  = "4 WHEEL STANDARD"
    $ws.Range("X" + $r).Value  = "STD"
    $ws.Range("Y" + $r).Value  = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
    $ws.Range("Z" + $r).Value  = 42
    $ws.Range("AA" + $r).Value = 42
    $ws.Range("AB" + $r).Value = "Y"
    $ws.Range("AC" + $r).Value = $rd.AC
    $ws.Range("AD" + $r).Value = $rd.AD
    $ws.Range("AE" + $r).Value = $rd.AE
    $ws.Range("AF" + $r).Value = $rd.AF
    $ws.Range("AG" + $r).Value = $rd.AG
    $ws.Range("AH" + $r).Value = $rd.AH
    $ws.Range("AI" + $r).Value = "Y"
    $ws.Range("AJ" + $r).Value = "N"

    # Copy cell formatting (left alignment + conditional Good/Bad highlight
    # on ENTRYDATE/VALID) from the matching existing row.
    $ws.Range("A" + $rd.SrcRow + ":AJ" + $rd.SrcRow).Copy()
    $ws.Range("A" + $r + ":AJ" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Match the final selection recorded in the sheet view
$ws.Range("E16").Select()
